# Add season record columns (Wins, Losses, Ties) to the PHI_2021 sheet.
# The original scraper only pulled team/player statistics and missed the
# season win-loss-tie record, so this adds it as three new trailing
# columns (AD, AE, AF).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold / centered / bordered header formatting used by the
# rest of row 1 (e.g. column A1) by copying its format onto the new
# header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-assert the header text/values (PasteSpecial only copied formats).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-57: season record (82-80-0) repeated for every player ---
$wins = 82
$losses = 80
$ties = 0

for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
